# Nexial command-list workbook update:
#   - add new "base" command: assertMatch(text,regex)
#   - add new "external" command: openFile(filePath)
#   - remove obsolete "tn.5250" target/category (and its column of commands)
#
# NOTE: this engine's Range.Insert/Delete with a vertical shift operates on
# the *whole row* (all columns), not just the addressed column, so the
# column-scoped list insertions below are done by writing the shifted
# values directly cell-by-cell instead of via Range.Insert/Delete.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) "base" command list (column F) - insert assertMatch(text,regex)
#    alphabetically before assertNotContain(text,substring); everything
#    from the old F11 onward shifts down one row (F11:F44 -> F12:F45).
# ---------------------------------------------------------------------
$ws.Range("F45").Value = "waitForCondition(conditions,maxWaitMs)"
$ws.Range("F44").Value = "waitFor(waitMs)"
$ws.Range("F43").Value = "verbose(text)"
$ws.Range("F42").Value = "substringBetween(text,start,end,saveVar)"
$ws.Range("F41").Value = "substringBefore(text,delim,saveVar)"
$ws.Range("F40").Value = "substringAfter(text,delim,saveVar)"
$ws.Range("F39").Value = "stopRecording()"
$ws.Range("F38").Value = "startRecording()"
$ws.Range("F37").Value = "split(text,delim,saveVar)"
$ws.Range("F36").Value = "section(steps)"
$ws.Range("F35").Value = "saveVariablesByRegex(var,regex)"
$ws.Range("F34").Value = "saveVariablesByPrefix(var,prefix)"
$ws.Range("F33").Value = "saveReplace(text,regex,replace,saveVar)"
$ws.Range("F32").Value = "saveMatches(text,regex,saveVar)"
$ws.Range("F31").Value = "saveCount(text,regex,saveVar)"
$ws.Range("F30").Value = "save(var,value)"
$ws.Range("F29").Value = "repeatUntil(steps,maxWaitMs)"
$ws.Range("F28").Value = "prependText(var,prependWith)"
$ws.Range("F27").Value = "outputToCloud(resource)"
$ws.Range("F26").Value = "macroFlex(macro,input,output)"
$ws.Range("F25").Value = "macro(file,sheet,name)"
$ws.Range("F24").Value = "incrementChar(var,amount,config)"
$ws.Range("F23").Value = "failImmediate(text)"
$ws.Range("F22").Value = "copyIntoClipboard(text)"
$ws.Range("F21").Value = "copyFromClipboard(var)"
$ws.Range("F20").Value = "clearClipboard()"
$ws.Range("F19").Value = "clear(vars)"
$ws.Range("F18").Value = "assertVarPresent(var)"
$ws.Range("F17").Value = "assertVarNotPresent(var)"
$ws.Range("F16").Value = "assertTextOrder(var,descending)"
$ws.Range("F15").Value = "assertStartsWith(text,prefix)"
$ws.Range("F14").Value = "assertNotEqual(expected,actual)"
$ws.Range("F13").Value = "assertNotEmpty(text)"
$ws.Range("F12").Value = "assertNotContain(text,substring)"
$ws.Range("F11").Value = "assertMatch(text,regex)"

# ---------------------------------------------------------------------
# 2) "external" command list (column J) - insert openFile(filePath)
#    alphabetically before runJUnit(className); everything from the old
#    J2 onward shifts down one row (J2:J6 -> J3:J7).
# ---------------------------------------------------------------------
$ws.Range("J7").Value = "terminate(programName)"
$ws.Range("J6").Value = "tail(id,file)"
$ws.Range("J5").Value = "runProgramNoWait(programPathAndParams)"
$ws.Range("J4").Value = "runProgram(programPathAndParams)"
$ws.Range("J3").Value = "runJUnit(className)"
$ws.Range("J2").Value = "openFile(filePath)"

# ---------------------------------------------------------------------
# 3) "target" list (column A) - remove "tn.5250" entry; everything below
#    it shifts up one row (A28:A33 -> A27:A32) and A33 becomes empty.
# ---------------------------------------------------------------------
$ws.Range("A27").Value = "web"
$ws.Range("A28").Value = "webalert"
$ws.Range("A29").Value = "webcookie"
$ws.Range("A30").Value = "ws"
$ws.Range("A31").Value = "ws.async"
$ws.Range("A32").Value = "xml"
$ws.Range("A33").ClearContents()

# ---------------------------------------------------------------------
# 4) Remove the obsolete "tn.5250" data column (AA) outright - this is a
#    genuine whole-column delete, so every column to its right (web,
#    webalert, webcookie, ws, ws.async, xml - AB..AG) shifts left one
#    column into AA..AF.
# ---------------------------------------------------------------------
$ws.Columns.Item(27).Delete()

# ---------------------------------------------------------------------
# 5) Update the defined names so they reference the new ranges.
# ---------------------------------------------------------------------
$wb.Names.Item("base").RefersTo = "='#system'!`$F`$2:`$F`$45"
$wb.Names.Item("external").RefersTo = "='#system'!`$J`$2:`$J`$7"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$32"
$wb.Names.Item("web").RefersTo = "='#system'!`$AA`$2:`$AA`$151"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AC`$2:`$AC`$10"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AD`$2:`$AD`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AE`$2:`$AE`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AF`$2:`$AF`$27"

Write-Output "done"
